$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (FECHA/HORA shift right to E/F)
$ws.Range("D:D").Insert()

# New header for the inserted column
$ws.Range("D1").Value = "TOKENS"

# New data cell under the inserted column (stored as text, not a number)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20"
$ws.Range("D2").Style = "Normal"

# Normalize name/area casing to lowercase
$ws.Range("B2").Value = "andres salcedo"
$ws.Range("C2").Value = "sistemas"

# Update the time value (now in column F after the insert)
$ws.Range("F2").Value = "16:04:15"
